{"js": "// Update the two-digit division worksheet: replace each \"N\u00f7N=\" expression\n// in the body (they live inside table cells) with a new expression, in\n// document order. Several old values reappear as new values elsewhere\n// (e.g. \"10\u00f74=\", \"77\u00f73=\", \"74\u00f78=\"), so replacement MUST be positional\n// (Nth occurrence -> Nth replacement) rather than a blind global\n// find/replace-by-text, otherwise a freshly written value could be\n// re-matched and overwritten by a later step.\nconst oldToNew = [\n  [\"67\u00f77=\", \"60\u00f79=\"],\n  [\"75\u00f77=\", \"10\u00f74=\"],\n  [\"12\u00f73=\", \"77\u00f73=\"],\n  [\"38\u00f76=\", \"62\u00f76=\"],\n  [\"19\u00f76=\", \"69\u00f73=\"],\n  [\"75\u00f75=\", \"68\u00f78=\"],\n  [\"71\u00f79=\", \"51\u00f75=\"],\n  [\"81\u00f75=\", \"85\u00f76=\"],\n  [\"19\u00f74=\", \"29\u00f73=\"],\n  [\"40\u00f77=\", \"28\u00f75=\"],\n  [\"52\u00f77=\", \"42\u00f78=\"],\n  [\"41\u00f77=\", \"92\u00f74=\"],\n  [\"57\u00f75=\", \"36\u00f79=\"],\n  [\"77\u00f73=\", \"64\u00f79=\"],\n  [\"71\u00f74=\", \"95\u00f75=\"],\n  [\"78\u00f73=\", \"89\u00f77=\"],\n  [\"74\u00f78=\", \"55\u00f79=\"],\n  [\"94\u00f74=\", \"46\u00f78=\"],\n  [\"10\u00f74=\", \"56\u00f73=\"],\n  [\"27\u00f76=\", \"68\u00f75=\"],\n  [\"88\u00f72=\", \"50\u00f74=\"],\n  [\"16\u00f76=\", \"52\u00f78=\"],\n  [\"85\u00f74=\", \"32\u00f76=\"],\n  [\"97\u00f72=\", \"74\u00f78=\"],\n  [\"46\u00f72=\", \"93\u00f72=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nconst divisionPattern = /^\\d+\u00f7\\d+=$/;\nlet matchIndex = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  if (!divisionPattern.test(text)) {\n    continue;\n  }\n  if (matchIndex >= oldToNew.length) {\n    break;\n  }\n  const [expectedOld, newText] = oldToNew[matchIndex];\n  if (text === expectedOld) {\n    const range = paragraph.getRange();\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  matchIndex++;\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit division worksheet: replace each \"N\u00f7N=\" expression\n# in the table with a new expression, in document order. Several old values\n# reappear as new values elsewhere (e.g. \"10\u00f74=\", \"77\u00f73=\", \"74\u00f78=\"), so the\n# replacement MUST be positional (Nth occurrence -> Nth replacement) rather\n# than a blind global find/replace-by-text; otherwise a freshly written\n# value could be re-matched and overwritten by a later step.\n$pairs = @(\n    @(\"67\u00f77=\", \"60\u00f79=\"),\n    @(\"75\u00f77=\", \"10\u00f74=\"),\n    @(\"12\u00f73=\", \"77\u00f73=\"),\n    @(\"38\u00f76=\", \"62\u00f76=\"),\n    @(\"19\u00f76=\", \"69\u00f73=\"),\n    @(\"75\u00f75=\", \"68\u00f78=\"),\n    @(\"71\u00f79=\", \"51\u00f75=\"),\n    @(\"81\u00f75=\", \"85\u00f76=\"),\n    @(\"19\u00f74=\", \"29\u00f73=\"),\n    @(\"40\u00f77=\", \"28\u00f75=\"),\n    @(\"52\u00f77=\", \"42\u00f78=\"),\n    @(\"41\u00f77=\", \"92\u00f74=\"),\n    @(\"57\u00f75=\", \"36\u00f79=\"),\n    @(\"77\u00f73=\", \"64\u00f79=\"),\n    @(\"71\u00f74=\", \"95\u00f75=\"),\n    @(\"78\u00f73=\", \"89\u00f77=\"),\n    @(\"74\u00f78=\", \"55\u00f79=\"),\n    @(\"94\u00f74=\", \"46\u00f78=\"),\n    @(\"10\u00f74=\", \"56\u00f73=\"),\n    @(\"27\u00f76=\", \"68\u00f75=\"),\n    @(\"88\u00f72=\", \"50\u00f74=\"),\n    @(\"16\u00f76=\", \"52\u00f78=\"),\n    @(\"85\u00f74=\", \"32\u00f76=\"),\n    @(\"97\u00f72=\", \"74\u00f78=\"),\n    @(\"46\u00f72=\", \"93\u00f72=\")\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$matchIndex = 0\nforeach ($cell in $table.Range.Cells) {\n    $raw = $cell.Range.Text\n    # Cell.Range.Text carries a trailing cell-mark (CR + BEL); strip it.\n    $text = $raw.Substring(0, $raw.Length - 2)\n    if ($text -match '^\\d+.\\d+=$') {\n        if ($matchIndex -lt $pairs.Count) {\n            $pair = $pairs[$matchIndex]\n            $expectedOld = $pair[0]\n            $newText = $pair[1]\n            if ($text -eq $expectedOld) {\n                $cell.Range.Text = $newText\n            }\n            $matchIndex++\n        }\n    }\n}\n"}
